$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New offense category introduced by this commit.
$offense = "Criminal Offenses - Forcible Sex Offenses"
$date = "sum2013"

$sectors = @(
  "Public, 4-year or above",
  "Private nonprofit, 4-year or above",
  "Private for-profit, 4-year or above",
  "Public, 2-year",
  "Private nonprofit, 2-year",
  "Private for-profit, 2-year",
  "Public, less-than 2-year",
  "Private nonprofit, less-than 2-year",
  "Private for-profit, less-than 2-year"
)

$locations = @(
  "On Campus (excluding Residence Halls)",
  "On Campus (Residence Halls)",
  "Non-Campus",
  "Public Property"
)

# counts[location] -> values aligned with $sectors order
$counts = @{
  "On Campus (excluding Residence Halls)" = @(615, 531, 5, 221, 12, 4, 10, 0, 3);
  "On Campus (Residence Halls)"           = @(1653, 1895, 15, 82, 3, 3, 1, 0, 3);
  "Non-Campus"                            = @(359, 162, 15, 45, 2, 4, 0, 0, 5);
  "Public Property"                       = @(98, 139, 18, 50, 8, 18, 4, 5, 13)
}

$lastRow = 541
$row = $lastRow + 1

foreach ($loc in $locations) {
  $vals = $counts[$loc]
  for ($i = 0; $i -lt $sectors.Count; $i++) {
    $ws.Cells.Item($row, 1).Value = $sectors[$i]
    $ws.Cells.Item($row, 2).Value = $loc
    $ws.Cells.Item($row, 3).Value = $offense
    $ws.Cells.Item($row, 4).Value = $date
    $ws.Cells.Item($row, 5).Value = $vals[$i]

    # Match the formatting pattern used by the rest of the table: columns
    # A, B, D and E carry the "quote-prefixed" style (s="1"), column C does
    # not. Replicate it by copying formats from the row directly above,
    # which is also shaped that way.
    $ws.Range("A" + ($row - 1) + ":E" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A" + $row + ":E" + $row).PasteSpecial(-4122) | Out-Null

    $row = $row + 1
  }
}

$excel.CutCopyMode = $false

$ws.Range("C578").Select() | Out-Null
